$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 35 ("Add new database location") moves to Complete, gets Started/Completed dates ---
$ws.Range("F35").Value = "Complete"

# G35/H35 previously blank; give them the same date-formatted style already used
# elsewhere in the sheet (e.g. G28/H28) instead of synthesizing a brand-new number
# format, then fill in the actual date values.
$ws.Range("G28:H28").Copy()
$ws.Range("G35:H35").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G35").Value = 43018   # 2017-10-10
$ws.Range("H35").Value = 43019   # 2017-10-11

# --- New row 37: "Ability to query ANY table from given DB" ---
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Ability to query ANY table from given DB"
$ws.Range("C37").Value = "HIGH"
$ws.Range("D37").Value = "LOW"
$ws.Range("E37").Value = "On page load, server asks DB for list of tables. Used to poplate window.availableTables."
$ws.Range("F37").Value = "In progress"

$ws.Range("G28").Copy()
$ws.Range("G37").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G37").Value = 43019   # 2017-10-11

$excel.CutCopyMode = 0

# --- Update the view so the new rows are visible, matching the saved view state ---
$ws.Range("G38").Select()
$excel.ActiveWindow.ScrollRow = 23
